# The document has two Pearson logo pictures (PearsonLogo.png) living in the
# two footers, and two BTEC logo pictures (BTec_Logo-Orange.jpg) living in the
# two headers. Each picture's drawing name ("image1.png"/"image2.jpg" etc.)
# needs to be swapped so that it matches the file name of the media part it
# is actually bound to:
#   footers: PearsonLogo picture  -> renamed from "image1.png" to "image2.png"
#   headers: BTec_Logo-Orange pic -> renamed from "image2.jpg" to "image1.jpg"
#
# An InlineShape does not expose a settable Name on its own (same as real
# Word's object model - only a converted/floating Shape does), but once the
# picture's own Range is selected, the InlineShape reachable from
# Selection.InlineShapes accepts a Name assignment and writes it back into
# the drawing's properties. We do this for every header and every footer so
# both logo pictures in each story get their corrected name.

$d = $word.ActiveDocument

function Rename-LogoPicture($headerFooter, [string]$newName) {
    if ($headerFooter.Exists -and $headerFooter.Range.InlineShapes.Count -ge 1) {
        $headerFooter.Range.InlineShapes.Item(1).Range.Select()
        $word.Selection.InlineShapes.Item(1).Name = $newName
    }
}

foreach ($sec in $d.Sections) {
    # Footers hold the Pearson logo: image1.png -> image2.png
    for ($fi = 1; $fi -le 3; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        Rename-LogoPicture $ftr "image2.png"
    }

    # Headers hold the BTEC logo: image2.jpg -> image1.jpg
    for ($hi = 1; $hi -le 3; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        Rename-LogoPicture $hdr "image1.jpg"
    }
}
